# Auto-generated edit script: updates column F ("想去人数") values
# across all 4 worksheets (展览, 演出, 本地生活, 全部类型) per the commit diff.
$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3549
$ws.Range("F5").Value = 8308
$ws.Range("F7").Value = 118
$ws.Range("F8").Value = 2217
$ws.Range("F10").Value = 92
$ws.Range("F11").Value = 69
$ws.Range("F12").Value = 638
$ws.Range("F13").Value = 103
$ws.Range("F14").Value = 7321
$ws.Range("F16").Value = 7607
$ws.Range("F18").Value = 57452
$ws.Range("F19").Value = 57452
$ws.Range("F20").Value = 4717
$ws.Range("F21").Value = 1055
$ws.Range("F22").Value = 931
$ws.Range("F23").Value = 493
$ws.Range("F25").Value = 922
$ws.Range("F28").Value = 4983
$ws.Range("F30").Value = 101
$ws.Range("F32").Value = 899
$ws.Range("F33").Value = 1333
$ws.Range("F34").Value = 1775
$ws.Range("F35").Value = 19
$ws.Range("F36").Value = 179
$ws.Range("F37").Value = 224
$ws.Range("F40").Value = 726
$ws.Range("F41").Value = 39
$ws.Range("F43").Value = 248
$ws.Range("F44").Value = 106
$ws.Range("F45").Value = 3
$ws.Range("F47").Value = 192
$ws.Range("F48").Value = 14
$ws.Range("F49").Value = 55

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 194
$ws.Range("F4").Value = 57
$ws.Range("F6").Value = 132
$ws.Range("F10").Value = 7589
$ws.Range("F11").Value = 123
$ws.Range("F13").Value = 5
$ws.Range("F20").Value = 71
$ws.Range("F42").Value = 120

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2355
$ws.Range("F5").Value = 1591
$ws.Range("F8").Value = 2406
$ws.Range("F9").Value = 9420
$ws.Range("F10").Value = 1746
$ws.Range("F11").Value = 178
$ws.Range("F12").Value = 111
$ws.Range("F15").Value = 260
$ws.Range("F16").Value = 2290
$ws.Range("F17").Value = 48
$ws.Range("F18").Value = 486

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2355
$ws.Range("F4").Value = 8308
$ws.Range("F6").Value = 1746
$ws.Range("F8").Value = 118
$ws.Range("F9").Value = 2290
$ws.Range("F10").Value = 69
$ws.Range("F11").Value = 7607
$ws.Range("F12").Value = 57452
$ws.Range("F13").Value = 194
$ws.Range("F15").Value = 57
$ws.Range("F16").Value = 4717
$ws.Range("F18").Value = 1055
$ws.Range("F19").Value = 931
$ws.Range("F21").Value = 922
$ws.Range("F23").Value = 4983
$ws.Range("F25").Value = 101
$ws.Range("F27").Value = 899
$ws.Range("F28").Value = 1333
$ws.Range("F29").Value = 1775
$ws.Range("F30").Value = 123
$ws.Range("F31").Value = 486
$ws.Range("F33").Value = 71
$ws.Range("F34").Value = 19
$ws.Range("F35").Value = 179
$ws.Range("F36").Value = 224
$ws.Range("F39").Value = 726
$ws.Range("F41").Value = 248
$ws.Range("F46").Value = 55
